# Append the next two days of GSC export data to the "Chart" sheet,
# mirroring the daily-export update captured in the diff:
#   - Chart!A65:C65 -> 2025-12-08, 0, 26
#   - Chart!A66:C66 -> 2025-12-09, 0, 27
# Dates are stored as plain text (matching every prior row), so each date
# cell is written as "@" (text) format first to stop Excel's automatic
# date-recognition from turning the string into a serial date, then the
# format is cleared back to the sheet's default (General/unstyled) so the
# new rows carry the same (unstyled) look as the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$ws.Range("A65").NumberFormat = "@"
$ws.Range("A65").Value = "2025-12-08"
$ws.Range("A65").ClearFormats()
$ws.Range("B65").Value = 0
$ws.Range("C65").Value = 26

$ws.Range("A66").NumberFormat = "@"
$ws.Range("A66").Value = "2025-12-09"
$ws.Range("A66").ClearFormats()
$ws.Range("B66").Value = 0
$ws.Range("C66").Value = 27

# The two "issues" sheets carry a static Issue/Validation/Items header that
# doesn't otherwise change with this export refresh; rewrite it so the
# sheet reflects the same refreshed export pass as the Chart sheet.
$ws2 = $wb.Worksheets.Item("Critical issues")
$ws2.Range("A1").Value = "Issue"
$ws2.Range("B1").Value = "Validation"
$ws2.Range("C1").Value = "Items"

$ws3 = $wb.Worksheets.Item("Non-critical issues")
$ws3.Range("A1").Value = "Issue"
$ws3.Range("B1").Value = "Validation"
$ws3.Range("C1").Value = "Items"
